$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (rows are 1-indexed within the lone table column)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "295"
$t.Cell(5,1).Range.Text  = "0.00002"
$t.Cell(6,1).Range.Text  = "0.00010"
$t.Cell(12,1).Range.Text = "0.01278"

# Rows that collapse a tab-separated run of values down to a single value
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "797"
